# Rename the picture shapes' "name" attribute (wp:docPr/@name and
# pic:cNvPr/@name) inside the document's headers/footers:
#   - Pearson logo shapes (id="1" and id="2"): image2.png -> image1.png
#   - BTec logo shape      (id="3"):           image1.jpg -> image2.jpg
#
# The relationship targets / actual media parts are untouched; this is a
# pure rename of the shape's display name, done by round-tripping each
# header/footer's WordOpenXML through a targeted regex substitution.

$d = $word.ActiveDocument

function Update-StoryXml($story) {
    if (-not $story.Exists) { return }

    $xml = $story.Range.WordOpenXML

    $updated = $xml
    $updated = $updated -replace '(descr="Y:\\Together Design\\Pearson Edexcel PowerPoint amends\\Assets\\PearsonLogo\.png" id="(?:1|2|0)" name=")image2\.png(")', '${1}image1.png${2}'
    $updated = $updated -replace '(descr="BTec_Logo-Orange" id="(?:3|0)" name=")image1\.jpg(")', '${1}image2.jpg${2}'

    if ($updated -ne $xml) {
        $story.Range.WordOpenXML = $updated
    }
}

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        Update-StoryXml $sec.Headers.Item($i)
        Update-StoryXml $sec.Footers.Item($i)
    }
}

Write-Host "Done"
